# CronogramaTSP.xlsx edit — "ECU-9 creado, Modificar informacion basica pacientes"
#
# Summary of the change (derived from the OOXML diff):
#  - Rename the three "ECU. Crear ..." tasks (rows 10-12) to "ECU. Registrar ...".
#  - Mark the "ECU. Registrar pacientes" / "ECU. Registrar propietarios" tasks
#    (rows 10-11) as finished: copy the "LISTO" row formatting (green fill,
#    column H = "LISTO") from an already-completed row, and clear the leftover
#    stray value in E10.
#  - Reassign the owner of "ECU. Registrar MVZ's" (row 12) from Jenny to Julián.
#  - The "ECU. Modificar información básica propietarios/pacientes" rows
#    (14-15) are now also finished, so they get the same "LISTO" formatting.
#  - Minor view bookkeeping: active cell / window position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 10 & 11: "ECU. Crear ..." -> "ECU. Registrar ...", now finished ---

# Copy the "done" look (green fill, border) from row 13 (already LISTO) onto
# rows 10 and 11 before touching their values.
$ws.Range("B13:I13").Copy() | Out-Null
$ws.Range("B10:I10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B13:I13").Copy() | Out-Null
$ws.Range("B11:I11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B10").Value = "ECU. Registrar pacientes"
$ws.Range("E10").ClearContents() | Out-Null
$ws.Range("H10").Value = "LISTO"

$ws.Range("B11").Value = "ECU. Registrar propietarios"
$ws.Range("H11").Value = "LISTO"

# --- Row 12: "ECU. Crear MVZ's" -> "ECU. Registrar MVZ's", owner Jenny -> Julián ---
$ws.Range("B12").Value = "ECU. Registrar MVZ's"
$ws.Range("I12").Value = "Julián"

# --- Rows 14 & 15: now finished too, pick up the LISTO formatting ---
$ws.Range("B13:I13").Copy() | Out-Null
$ws.Range("B14:I14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B13:I13").Copy() | Out-Null
$ws.Range("B15:I15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H14").Value = "LISTO"
$ws.Range("H15").Value = "LISTO"

# --- View bookkeeping ---
$ws.Range("D17").Select() | Out-Null
